$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Trim the workbook down to a single sheet and rename it Sheet1
#    (index-based access avoids any cyrillic-name round-trip issues)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$wb.Worksheets.Item(3).Delete()
$wb.Worksheets.Item(2).Delete()
$ws.Name = "Sheet1"

# ---------------------------------------------------------------------------
# 2. Insert a new "properties" row right above the existing C6:E6 header row
#    (everything from row 6 down shifts by one row)
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Insert()

$newRow = $ws.Range("C6:E6")
$newRow.ClearFormats()
$newRow.Borders.LineStyle = 1
$newRow.Borders.Weight = 2

$ws.Range("C6").Value = "properties"
$ws.Range("D6").Value = "lob"
$ws.Range("E6").Value = "lob1"

# ---------------------------------------------------------------------------
# 3. Update the saved selection to match the authored view
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E18").Select()

Write-Host "done"
